$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68").Value = "Mortarpod"
$ws.Range("A46").Value = "Dockside Extortionist"
$ws.Range("A29").Value = "Lena, Selfless Champion"

$ws.Range("A29").Select()
$excel.ActiveWindow.ScrollRow = 16
